$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 8500
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 7000
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = -9707
$ws.Range("N10").Value = -7586
$ws.Range("H11").Value = 702.3333
$ws.Range("I11").Value = 702.3333
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 702.3333
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -562.3333
$ws.Range("H40").Value = 46166.11
$ws.Range("I40").Value = 100961.25
$ws.Range("J40").Value = 2330
$ws.Range("K40").Value = 100961.25
$ws.Range("L40").Value = 2330
$ws.Range("M40").Value = -100786.25
$ws.Range("N40").Value = -2680
$ws.Range("H69").Value = 3360.75
$ws.Range("I69").Value = 3213
$ws.Range("J69").Value = 3410
$ws.Range("K69").Value = 9639
$ws.Range("L69").Value = 10230
$ws.Range("M69").Value = -8765
$ws.Range("N69").Value = -11978
$ws.Range("H72").Value = 3360.75
$ws.Range("I72").Value = 3213
$ws.Range("J72").Value = 3410
$ws.Range("K72").Value = 28917
$ws.Range("L72").Value = 30690
$ws.Range("M72").Value = -24549
$ws.Range("N72").Value = -39426
$ws.Range("H108").Value = 32000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 32000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H137").Value = 4550150
$ws.Range("I137").Value = 9097391
$ws.Range("J137").Value = 2909.3635
$ws.Range("K137").Value = 27292173
$ws.Range("L137").Value = 8728.0905
$ws.Range("M137").Value = -27289623
$ws.Range("H138").Value = 3317.3193
$ws.Range("I138").Value = 1630.0312
$ws.Range("J138").Value = 4667.15
$ws.Range("K138").Value = 4890.0936
$ws.Range("L138").Value = 14001.45
$ws.Range("M138").Value = 249.9063999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7570.14
$ws.Range("I32").Value = 5190.643
$ws.Range("J32").Value = 20062.5
$ws.Range("K32").Value = 5190.643
$ws.Range("L32").Value = 20062.5
$ws.Range("M32").Value = -4903.643
$ws.Range("H102").Value = 2582.963
$ws.Range("I102").Value = 2230.476
$ws.Range("J102").Value = 3816.6667
$ws.Range("K102").Value = 2230.476
$ws.Range("L102").Value = 3816.6667
$ws.Range("M102").Value = -608.4760000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2276.25
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 3001.6667
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 3001.6667
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = -3281.6667
$ws.Range("H86").Value = 1907.9584
$ws.Range("I86").Value = 1266.1666
$ws.Range("J86").Value = 3833.3333
$ws.Range("K86").Value = 1266.1666
$ws.Range("L86").Value = 3833.3333
$ws.Range("M86").Value = -143.1666
$ws.Range("N86").Value = -6079.3333
$ws.Range("H89").Value = 1907.9584
$ws.Range("I89").Value = 1266.1666
$ws.Range("J89").Value = 3833.3333
$ws.Range("K89").Value = 6330.833000000001
$ws.Range("L89").Value = 19166.6665
$ws.Range("M89").Value = -714.8330000000005
$ws.Range("N89").Value = -30398.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 13132.909
$ws.Range("I10").Value = 494
$ws.Range("J10").Value = 70008
$ws.Range("K10").Value = 494
$ws.Range("L10").Value = 70008
$ws.Range("M10").Value = -355
$ws.Range("N10").Value = -70286
$ws.Range("H31").Value = 3231934.2
$ws.Range("I31").Value = 9093914
$ws.Range("J31").Value = 7845.65
$ws.Range("K31").Value = 9093914
$ws.Range("L31").Value = 7845.65
$ws.Range("M31").Value = -9093619
$ws.Range("N31").Value = -8435.65
$ws.Range("H34").Value = 3231934.2
$ws.Range("I34").Value = 9093914
$ws.Range("J34").Value = 7845.65
$ws.Range("K34").Value = 9093914
$ws.Range("L34").Value = 7845.65
$ws.Range("M34").Value = -9093712
$ws.Range("N34").Value = -8249.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1668.04
$ws.Range("I9").Value = 447.2857
$ws.Range("J9").Value = 2142.7778
$ws.Range("K9").Value = 1341.8571
$ws.Range("L9").Value = 6428.3334
$ws.Range("M9").Value = -1117.8571
$ws.Range("N9").Value = -6876.3334
$ws.Range("H101").Value = 9343
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 9343
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 28029
$ws.Range("N101").Value = -32897
$ws.Range("H122").Value = 1408.5555
$ws.Range("I122").Value = 490
$ws.Range("J122").Value = 1671
$ws.Range("K122").Value = 4410
$ws.Range("L122").Value = 15039
$ws.Range("M122").Value = -1960
$ws.Range("N122").Value = -19939
$ws.Range("H131").Value = 1500.575
$ws.Range("I131").Value = 967.7273
$ws.Range("J131").Value = 2151.8333
$ws.Range("K131").Value = 2903.1819
$ws.Range("L131").Value = 6455.499899999999
$ws.Range("M131").Value = 2136.8181
$ws.Range("N131").Value = -16535.4999
$ws.Range("H132").Value = 2022.9642
$ws.Range("I132").Value = 1409.3846
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 12684.4614
$ws.Range("L132").Value = 89995.5
$ws.Range("M132").Value = -10154.4614
$ws.Range("N132").Value = -95055.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 22831.428
$ws.Range("I24").Value = 5806
$ws.Range("J24").Value = 25669
$ws.Range("K24").Value = 5806
$ws.Range("L24").Value = 25669
$ws.Range("M24").Value = -5633
$ws.Range("N24").Value = -26015
$ws.Range("H42").Value = 37932.4
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 37932.4
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 37932.4
$ws.Range("N42").Value = -38902.4
$ws.Range("H115").Value = 37932.4
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 37932.4
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 37932.4
$ws.Range("N115").Value = -40282.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1834.1875
$ws.Range("I93").Value = 1303.9
$ws.Range("J93").Value = 2718
$ws.Range("K93").Value = 1303.9
$ws.Range("L93").Value = 2718
$ws.Range("M93").Value = -55.90000000000009
$ws.Range("N93").Value = -5214
$ws.Range("H132").Value = 3118.5334
$ws.Range("I132").Value = 2147.5
$ws.Range("J132").Value = 4228.2856
$ws.Range("K132").Value = 6442.5
$ws.Range("L132").Value = 12684.8568
$ws.Range("M132").Value = -3912.5
$ws.Range("N132").Value = -17744.8568
$ws.Range("H136").Value = 3706383.2
$ws.Range("I136").Value = 5265155.5
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 15795466.5
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -15792916.5
$ws.Range("N136").Value = -18000
$ws.Range("H140").Value = 40000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 522.9231
$ws.Range("I100").Value = 483.16666
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 966.33332
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -425.33332
$ws.Range("H132").Value = 725715.1
$ws.Range("I132").Value = 1432144
$ws.Range("J132").Value = 19286.285
$ws.Range("K132").Value = 4296432
$ws.Range("L132").Value = 57858.855
$ws.Range("M132").Value = -4293902
$ws.Range("N132").Value = -62918.855
